# Update Name of Algo
# Applies updated KNN-imputed values to the result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.64
$ws.Range("E6").Value = 16.398
$ws.Range("A7").Value = -21.304
$ws.Range("B7").Value = 5.872
$ws.Range("B15").Value = 5.042
$ws.Range("A16").Value = -21.918
$ws.Range("C16").Value = -13.148
$ws.Range("C19").Value = -12.259
$ws.Range("E19").Value = 16.664
$ws.Range("E20").Value = 16.326
$ws.Range("B21").Value = 8.324999999999999
$ws.Range("E21").Value = 16.655
$ws.Range("B22").Value = 6.710000000000001
$ws.Range("B23").Value = 7.083999999999999
$ws.Range("E24").Value = 17.143
$ws.Range("A28").Value = -21.767
$ws.Range("A29").Value = -21.533
$ws.Range("A32").Value = -21.684
$ws.Range("B34").Value = 7.821
$ws.Range("D34").Value = -7.929
$ws.Range("E35").Value = 16.4
$ws.Range("C36").Value = -12.669
$ws.Range("E39").Value = 16.573
$ws.Range("A40").Value = -20.444
$ws.Range("E41").Value = 16.489
$ws.Range("B43").Value = 5.922000000000001
$ws.Range("D43").Value = -8.118
$ws.Range("B45").Value = 5.672
$ws.Range("C46").Value = -13.805
$ws.Range("D48").Value = -7.702
$ws.Range("B50").Value = 5.715000000000001
$ws.Range("C50").Value = -13.262
$ws.Range("B51").Value = 6.284
$ws.Range("A52").Value = -21.712
$ws.Range("A57").Value = -22.32
$ws.Range("A66").Value = -21.398
$ws.Range("B66").Value = 5.859
$ws.Range("B67").Value = 5.571
$ws.Range("D70").Value = -6.936
$ws.Range("D73").Value = -8.270999999999999
$ws.Range("E73").Value = 16.402
$ws.Range("E76").Value = 16.578
$ws.Range("B79").Value = 5.495
$ws.Range("B84").Value = 5.453
$ws.Range("E85").Value = 16.794
$ws.Range("D87").Value = -7.933
$ws.Range("B92").Value = 5.673
$ws.Range("D92").Value = -6.328
$ws.Range("C95").Value = -12.25
$ws.Range("B97").Value = 6.433
$ws.Range("C97").Value = -13.439
$ws.Range("E98").Value = 16.518
$ws.Range("A100").Value = -21.907
$ws.Range("D101").Value = -8.051
